# M11 Froze Token Embeddings + Decoder 123
# Update per-epoch accuracy values in column B of the Epoch Accuracy sheet
# and restore the selection to the data range (A2:B116).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B3"   = 0.84375
    "B4"   = 0.765625
    "B5"   = 0.609375
    "B6"   = 0.625
    "B8"   = 0.59375
    "B9"   = 0.59375
    "B10"  = 0.578125
    "B11"  = 0.578125
    "B12"  = 0.53125
    "B14"  = 0.546875
    "B15"  = 0.5
    "B16"  = 0.546875
    "B17"  = 0.5625
    "B18"  = 0.546875
    "B19"  = 0.5625
    "B20"  = 0.546875
    "B21"  = 0.546875
    "B22"  = 0.546875
    "B23"  = 0.578125
    "B28"  = 0.5625
    "B29"  = 0.5625
    "B30"  = 0.5625
    "B31"  = 0.5625
    "B32"  = 0.5625
    "B33"  = 0.5625
    "B34"  = 0.5625
    "B35"  = 0.5625
    "B36"  = 0.5625
    "B37"  = 0.5625
    "B38"  = 0.5625
    "B45"  = 0.546875
    "B46"  = 0.546875
    "B47"  = 0.546875
    "B48"  = 0.546875
    "B49"  = 0.546875
    "B50"  = 0.546875
    "B51"  = 0.546875
    "B52"  = 0.546875
    "B53"  = 0.546875
    "B54"  = 0.546875
    "B55"  = 0.546875
    "B56"  = 0.546875
    "B57"  = 0.546875
    "B58"  = 0.546875
    "B59"  = 0.546875
    "B60"  = 0.546875
    "B61"  = 0.546875
    "B62"  = 0.546875
    "B63"  = 0.546875
    "B64"  = 0.546875
    "B65"  = 0.546875
    "B66"  = 0.546875
    "B67"  = 0.546875
    "B68"  = 0.546875
    "B69"  = 0.546875
    "B70"  = 0.546875
    "B71"  = 0.546875
    "B72"  = 0.546875
    "B73"  = 0.546875
    "B74"  = 0.546875
    "B75"  = 0.546875
    "B76"  = 0.546875
    "B77"  = 0.546875
    "B78"  = 0.546875
    "B79"  = 0.546875
    "B80"  = 0.546875
    "B81"  = 0.546875
    "B82"  = 0.546875
    "B83"  = 0.546875
    "B84"  = 0.546875
    "B85"  = 0.546875
    "B86"  = 0.546875
    "B87"  = 0.546875
    "B88"  = 0.546875
    "B89"  = 0.546875
    "B90"  = 0.546875
    "B91"  = 0.546875
    "B92"  = 0.546875
    "B93"  = 0.546875
    "B94"  = 0.546875
    "B95"  = 0.546875
    "B96"  = 0.546875
    "B97"  = 0.546875
    "B98"  = 0.546875
    "B99"  = 0.546875
    "B100" = 0.546875
    "B101" = 0.546875
    "B102" = 0.546875
    "B103" = 0.484375
    "B104" = 0.546875
    "B105" = 0.609375
    "B106" = 0.5625
    "B107" = 0.59375
    "B108" = 0.59375
    "B110" = 0.6875
    "B111" = 0.46875
    "B112" = 0.453125
    "B113" = 0.640625
    "B114" = 0.546875
    "B115" = 0.515625
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# Select the full data range, matching the saved selection in the workbook.
$ws.Range("A2:B116").Select()
